# Update Name of Algo
# Apply updated RandomForest imputation results to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C7"   = -13.515
    "B9"   = 6.374599999999997
    "C12"  = -11.1801
    "D15"  = -8.774999999999999
    "B18"  = 6.760399999999999
    "B20"  = 9.374799999999995
    "C26"  = -12.93930000000001
    "B27"  = 5.675500000000002
    "C27"  = -13.1176
    "C29"  = -11.3435
    "C37"  = -13.8049
    "C38"  = -13.717
    "D38"  = -8.664499999999991
    "D44"  = -7.221600000000003
    "C51"  = -12.29069999999999
    "D51"  = -7.772000000000003
    "C55"  = -14.07680000000001
    "D57"  = -8.044400000000001
    "D63"  = -7.8311
    "B69"  = 6.389099999999997
    "C69"  = -12.7691
    "C70"  = -12.11649999999999
    "D70"  = -7.963799999999996
    "B76"  = 4.8313
    "B82"  = 7.060900000000003
    "C83"  = -14.1983
    "D99"  = -7.6495
    "C102" = -13.3159
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
